$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '61.693.65'
$ws.Cells.Item(2, 5).Value = '  +0.32%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.448.08'
$ws.Cells.Item(3, 5).Value = '  +2.56%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'578.07"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.93%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'147.24"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +7.94%  '

# Row 7
$ws.Cells.Item(7, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(7, 4).Value = '3.448.96'
$ws.Cells.Item(7, 5).Value = '  +2.63%  '

# Row 8
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.473"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.67%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'7.67"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +2.49%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.124"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.23%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.386"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -0.60%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '4.032.35'
$ws.Cells.Item(13, 5).Value = '  +2.45%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'27.91"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +8.01%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.53%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'0.0000175"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +0.42%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.445.40'
$ws.Cells.Item(17, 5).Value = '  +2.43%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '61.699.42'
$ws.Cells.Item(18, 5).Value = '  +0.18%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'6.27"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +7.25%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'14.06"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.74%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'9.37"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.37%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'382.88"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.10%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'0.565"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +2.13%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '3.586.51'
$ws.Cells.Item(24, 5).Value = '  +2.46%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.03%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.38%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'72.22"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +1.36%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.0000124"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -1.23%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'0.179"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +8.71%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'7.71"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +2.30%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'1.56"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -12.72%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.09%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.98%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'24.00"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.90%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.19%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.02%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +2.02%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'165.73"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.03%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.0782"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +2.65%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = "'0.791"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +2.90%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = "'25.56"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +5.53%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.10%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'1.73"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.86%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'OKB'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(46, 4).Value = "'42.10"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +1.46%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Filecoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(47, 4).Value = "'4.47"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +1.81%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).Value = "'1.17"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -2.74%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '2.606.72'
$ws.Cells.Item(49, 5).Value = '  +9.93%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'23.51"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.67%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'6.85"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.34%  '
